$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.747.99'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.603.44'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'211.96"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = "'0.248"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').Value = "'19.74"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('D11').Value = "'0.0847"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('D12').Value = '1.827.81'
$ws.Range('D13').Value = '1.591.14'
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = "'0.523"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = "'65.15"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = "'210.55"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('E19').Value = '  +2.29%  '
$ws.Range('D21').Value = "'4.28"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('E22').Value = '  -2.32%  '
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').Value = "'143.69"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').Value = "'7.11"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = "'15.40"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('D33').Value = '1.295.54'
$ws.Range('E33').Value = '  +0.76%  '
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('E36').Value = '  -3.15%  '
$ws.Range('D37').Value = "'1.18"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.35%  '
$ws.Range('D38').Value = "'0.0170"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  -2.16%  '
$ws.Range('D41').Value = "'2.20"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = "'0.787"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').Value = "'62.96"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').Value = '1.739.30'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').Value = "'90.70"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = "'1.56"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = "'0.102"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.0517"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'7.46"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = "'0.396"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.83%  '
